$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (changed date) column for rows 2-5 from 2023-10-13 (45212) to 2023-10-22 (45221)
$ws.Range("C2:C5").Value = 45221
